# NIT-9005241478.xlsx — "Actualiza base de datos EC y agrega parte 1 de
# nuevos estado de cuenta"
#
# The statement currently lists two workers (BIENVENIDO DE LOS RIOS
# SANMARTIN and, across two overdue periods, DAIRA MARIA NORIEGA TEHERAN).
# This edit refreshes the data for "parte 1" of the new account statement:
# the second worker's two detail rows are removed, the per-worker salary
# figure for the remaining worker is updated, and the summary totals
# (overdue value, worker count, period count) are refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the detail rows (17 & 18) belonging to DAIRA MARIA NORIEGA
# TEHERAN; Excel shifts the rows below (the signature block) up by two.
$ws.Rows("17:18").Delete()

# Refresh the summary header figures.
$ws.Range("E11").Value = 33160   # VALOR MORA
$ws.Range("C13").Value = 1       # Cant. Trabajadores
$ws.Range("F13").Value = 1       # Cant. Periodos

# Update the remaining worker's Salario Basico.
$ws.Range("G16").Value = 1855117
